$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update res_bus vm_pu values: slack bus voltage set-point changed from 1.05 to 1.02 pu
# (case with 380 kV), and downstream bus voltages recalculated accordingly.

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.021200608447104
$ws.Range("D2").Value = 1.031702065088858
$ws.Range("E2").Value = 1.022107928532126
$ws.Range("F2").Value = 1.040270704415904
$ws.Range("I2").Value = 1.030121902728916
$ws.Range("J2").Value = 1.026393415012565
$ws.Range("K2").Value = 1.034509420442672
$ws.Range("L2").Value = 1.024943282708077
$ws.Range("M2").Value = 1.043053527838588
$ws.Range("N2").Value = 1.012856759222974

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.022071785405448
$ws.Range("D3").Value = 1.032392839003536
$ws.Range("E3").Value = 1.022845062441649
$ws.Range("F3").Value = 1.041196600686422
$ws.Range("I3").Value = 1.030242815009767
$ws.Range("J3").Value = 1.026902516367558
$ws.Range("K3").Value = 1.035009129008677
$ws.Range("L3").Value = 1.025487202339905
$ws.Range("M3").Value = 1.043789509210371
$ws.Range("N3").Value = 1.013027125019996

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022636124202438
$ws.Range("D4").Value = 1.032840145503898
$ws.Range("E4").Value = 1.023322971393596
$ws.Range("F4").Value = 1.041796577585247
$ws.Range("I4").Value = 1.030319790607528
$ws.Range("J4").Value = 1.027231945731804
$ws.Range("K4").Value = 1.035332127502101
$ws.Range("L4").Value = 1.025839423472659
$ws.Range("M4").Value = 1.044265971177639
$ws.Range("N4").Value = 1.013137310523711

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.022873521425757
$ws.Range("D5").Value = 1.033028270347061
$ws.Range("E5").Value = 1.023524106252417
$ws.Range("F5").Value = 1.042049012026028
$ws.Range("I5").Value = 1.030351848207693
$ws.Range("J5").Value = 1.027370438280265
$ws.Range("K5").Value = 1.035467831839639
$ws.Range("L5").Value = 1.025987560310855
$ws.Range("M5").Value = 1.044466329954137
$ws.Range("N5").Value = 1.013183619484734

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022913390172193
$ws.Range("D6").Value = 1.033059861812655
$ws.Range("E6").Value = 1.023557890646033
$ws.Range("F6").Value = 1.042091408799953
$ws.Range("I6").Value = 1.030357213036608
$ws.Range("J6").Value = 1.027393691780527
$ws.Range("K6").Value = 1.035490612217742
$ws.Range("L6").Value = 1.026012436799645
$ws.Range("M6").Value = 1.044499974205028
$ws.Range("N6").Value = 1.013191394186425

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022639295731279
$ws.Range("D7").Value = 1.03284265893611
$ws.Range("E7").Value = 1.023325658097471
$ws.Range("F7").Value = 1.041799949826623
$ws.Range("I7").Value = 1.030320220153901
$ws.Range("J7").Value = 1.027233796275741
$ws.Range("K7").Value = 1.035333941122708
$ws.Range("L7").Value = 1.025841402638097
$ws.Range("M7").Value = 1.044268648170553
$ws.Range("N7").Value = 1.013137929358049

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02149489568521
$ws.Range("D8").Value = 1.031935446094655
$ws.Range("E8").Value = 1.022356851988221
$ws.Range("F8").Value = 1.040583436249551
$ws.Range("I8").Value = 1.030163026672404
$ws.Range("J8").Value = 1.02656546609838
$ws.Range("K8").Value = 1.034678370221995
$ws.Range("L8").Value = 1.025127046483624
$ws.Range("M8").Value = 1.043302207041432
$ws.Range("N8").Value = 1.012914345758751

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019483196456033
$ws.Range("D9").Value = 1.030339421074689
$ws.Range("E9").Value = 1.020656919181664
$ws.Range("F9").Value = 1.038446448903213
$ws.Range("I9").Value = 1.029876388955239
$ws.Range("J9").Value = 1.025387886340053
$ws.Range("K9").Value = 1.03352057563722
$ws.Range("L9").Value = 1.023870374064359
$ws.Range("M9").Value = 1.041601066649087
$ws.Range("N9").Value = 1.012519979870855

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018145427760244
$ws.Range("D10").Value = 1.029277254327686
$ws.Range("E10").Value = 1.01952858940274
$ws.Range("F10").Value = 1.037026376199088
$ws.Range("I10").Value = 1.029678852002009
$ws.Range("J10").Value = 1.024602972252638
$ws.Range("K10").Value = 1.032747045700967
$ws.Range("L10").Value = 1.023034089433672
$ws.Range("M10").Value = 1.040468311655979
$ws.Range("N10").Value = 1.012256838013584

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.017566974435679
$ws.Range("D11").Value = 1.028817784630482
$ws.Range("E11").Value = 1.019041207280538
$ws.Range("F11").Value = 1.036412578556918
$ws.Range("I11").Value = 1.029591795603021
$ws.Range("J11").Value = 1.024263143378655
$ws.Range("K11").Value = 1.032411719675393
$ws.Range("L11").Value = 1.022672339875697
$ws.Range("M11").Value = 1.039978153288451
$ws.Range("N11").Value = 1.012142845348041

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017352234087929
$ws.Range("D12").Value = 1.028647187288306
$ws.Range("E12").Value = 1.018860352659051
$ws.Range("F12").Value = 1.036184754281521
$ws.Range("I12").Value = 1.029559231000105
$ws.Range("J12").Value = 1.024136923596144
$ws.Range("K12").Value = 1.032287108501393
$ws.Range("L12").Value = 1.02253802642463
$ws.Range("M12").Value = 1.039796138090165
$ws.Range("N12").Value = 1.012100496249709

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017398291034552
$ws.Range("D13").Value = 1.02868377779468
$ws.Range("E13").Value = 1.018899138380501
$ws.Range("F13").Value = 1.036233615754145
$ws.Range("I13").Value = 1.029566226524556
$ws.Range("J13").Value = 1.024163997796624
$ws.Range("K13").Value = 1.032313840541565
$ws.Range("L13").Value = 1.022566834541282
$ws.Range("M13").Value = 1.039835178615046
$ws.Range("N13").Value = 1.012109580591391

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.017549221414962
$ws.Range("D14").Value = 1.028803681559719
$ws.Range("E14").Value = 1.019026254083047
$ws.Range("F14").Value = 1.036393743104184
$ws.Range("I14").Value = 1.029589108449908
$ws.Range("J14").Value = 1.024252709846357
$ws.Range("K14").Value = 1.032401420414946
$ws.Range("L14").Value = 1.022661236325178
$ws.Range("M14").Value = 1.03996310679691
$ws.Range("N14").Value = 1.012139344901606

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.017642230939996
$ws.Range("D15").Value = 1.028877567575312
$ws.Range("E15").Value = 1.019104598271142
$ws.Range("F15").Value = 1.03649242509703
$ws.Range("I15").Value = 1.029603176568365
$ws.Range("J15").Value = 1.024307369335478
$ws.Range("K15").Value = 1.032455373862948
$ws.Range("L15").Value = 1.02271940789703
$ws.Range("M15").Value = 1.040041934435435
$ws.Range("N15").Value = 1.012157682741636

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018183834971809
$ws.Range("D16").Value = 1.029307757562737
$ws.Range("E16").Value = 1.019560960610496
$ws.Range("F16").Value = 1.037067135320108
$ws.Range("I16").Value = 1.02968459764334
$ws.Range("J16").Value = 1.024625526604567
$ws.Range("K16").Value = 1.032769292238126
$ws.Range("L16").Value = 1.023058105405733
$ws.Range("M16").Value = 1.040500848981076
$ws.Range("N16").Value = 1.012264402304939

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.018523786528517
$ws.Range("D17").Value = 1.029577727548949
$ws.Range("E17").Value = 1.019847545011463
$ws.Range("F17").Value = 1.03742793245739
$ws.Range("I17").Value = 1.029735263997444
$ws.Range("J17").Value = 1.024825110885115
$ws.Range("K17").Value = 1.032966103465707
$ws.Range("L17").Value = 1.023270660615036
$ws.Range("M17").Value = 1.040788803944356
$ws.Range("N17").Value = 1.012331331457328

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.018722152517857
$ws.Range("D18").Value = 1.029735240309199
$ws.Range("E18").Value = 1.020014819752291
$ws.Range("F18").Value = 1.037638485649437
$ws.Range("I18").Value = 1.029764669970983
$ws.Range("J18").Value = 1.02494152919661
$ws.Range("K18").Value = 1.033080863135925
$ws.Range("L18").Value = 1.023394675817153
$ws.Range("M18").Value = 1.040956795132923
$ws.Range("N18").Value = 1.012370365197025

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018789803381933
$ws.Range("D19").Value = 1.029788955472576
$ws.Range("E19").Value = 1.020071875554942
$ws.Range("F19").Value = 1.037710296828165
$ws.Range("I19").Value = 1.029774671720205
$ws.Range("J19").Value = 1.024981225505397
$ws.Range("K19").Value = 1.033119986904204
$ws.Range("L19").Value = 1.023436967752564
$ws.Range("M19").Value = 1.041014081155383
$ws.Range("N19").Value = 1.012383673848563

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.018487304866118
$ws.Range("D20").Value = 1.029548757777895
$ws.Range("E20").Value = 1.019816785319683
$ws.Range("F20").Value = 1.037389211328636
$ws.Range("I20").Value = 1.029729843162415
$ws.Range("J20").Value = 1.024803696940154
$ws.Range("K20").Value = 1.032944991300133
$ws.Range("L20").Value = 1.02324785178762
$ws.Range("M20").Value = 1.040757905787992
$ws.Range("N20").Value = 1.012324151096595

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.017504772774045
$ws.Range("D21").Value = 1.028768370938795
$ws.Range("E21").Value = 1.018988816672352
$ws.Range("F21").Value = 1.036346584952805
$ws.Range("I21").Value = 1.029582376581575
$ws.Range("J21").Value = 1.024226586134335
$ws.Range("K21").Value = 1.032375631864875
$ws.Range("L21").Value = 1.022633435784434
$ws.Range("M21").Value = 1.039925433703233
$ws.Range("N21").Value = 1.012130580247852

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.016887728017595
$ws.Range("D22").Value = 1.028278117543804
$ws.Range("E22").Value = 1.018469286810193
$ws.Range("F22").Value = 1.035692014330831
$ws.Range("I22").Value = 1.029488339775138
$ws.Range("J22").Value = 1.023863779325555
$ws.Range("K22").Value = 1.032017328723467
$ws.Range("L22").Value = 1.022247455241014
$ws.Range("M22").Value = 1.039402323873921
$ws.Range("N22").Value = 1.01200883344113

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017214767151786
$ws.Range("D23").Value = 1.028537971010254
$ws.Range("E23").Value = 1.018744599620767
$ws.Range("F23").Value = 1.03603892207534
$ws.Range("I23").Value = 1.029538315256313
$ws.Range("J23").Value = 1.024056105362029
$ws.Range("K23").Value = 1.03220730229996
$ws.Range("L23").Value = 1.022452039391852
$ws.Range("M23").Value = 1.039679605484442
$ws.Range("N23").Value = 1.012073377478526

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.018503789118165
$ws.Range("D24").Value = 1.02956184783521
$ws.Range("E24").Value = 1.019830683945526
$ws.Range("F24").Value = 1.037406707411656
$ws.Range("I24").Value = 1.02973229305825
$ws.Range("J24").Value = 1.024813372967355
$ws.Range("K24").Value = 1.032954531093431
$ws.Range("L24").Value = 1.023258158005877
$ws.Range("M24").Value = 1.040771867236372
$ws.Range("N24").Value = 1.012327395607735

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020002682080956
$ws.Range("D25").Value = 1.030751713189443
$ws.Range("E25").Value = 1.021095525823917
$ws.Range("F25").Value = 1.038998110770267
$ws.Range("I25").Value = 1.029951630467125
$ws.Range("J25").Value = 1.025692299512759
$ws.Range("K25").Value = 1.033820192782446
$ws.Range("L25").Value = 1.024194995670774
$ws.Range("M25").Value = 1.042040623104059
$ws.Range("N25").Value = 1.012621975588825

Write-Host "Updated vm_pu values for rows 2-25 (380 kV case)"